$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 4.277776333333333
$ws.Range("H2").Value = 12.833329
$ws.Range("I2").Value = 0.3536657835996513
$ws.Range("J2").Value = 0.3568846407551645
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.2791963333333333
$ws.Range("N2").Value = 0.837589
$ws.Range("O2").Value = 0.008912184157424861
$ws.Range("P2").Value = 0.009090659364840875
$ws.Range("Q2").Value = 1.194339467086778
$ws.Range("R2").Value = 10.749055203781
$ws.Range("S2").Value = 0.003151934593620062
$ws.Range("T2").Value = 0.003244316701648807

$ws.Range("G3").Value = 4.277776333333333
$ws.Range("H3").Value = 12.833329
$ws.Range("I3").Value = 0.3536657835996513
$ws.Range("J3").Value = 0.3568846407551645
$ws.Range("M3").Value = 2.510701
$ws.Range("N3").Value = 7.532103
$ws.Range("O3").Value = 0.08014370894160773
$ws.Range("P3").Value = 0.0817486651256118
$ws.Range("Q3").Value = 10.74021731787633
$ws.Range("R3").Value = 96.66195586088699
$ws.Range("S3").Value = 0.02834408762341608
$ws.Range("T3").Value = 0.02917484298556821

$ws.Range("G4").Value = 4.277776333333333
$ws.Range("H4").Value = 12.833329
$ws.Range("I4").Value = 0.3536657835996513
$ws.Range("J4").Value = 0.3568846407551645
$ws.Range("M4").Value = 16.644438
$ws.Range("N4").Value = 49.933314
$ws.Range("O4").Value = 0.5313046016107197
$ws.Range("P4").Value = 0.5419444960853593
$ws.Range("Q4").Value = 71.201182958034
$ws.Range("R4").Value = 640.810646622306
$ws.Range("S4").Value = 0.1879042582587558
$ws.Range("T4").Value = 0.1934116667946621

$ws.Range("G5").Value = 4.277776333333333
$ws.Range("H5").Value = 12.833329
$ws.Range("I5").Value = 0.3536657835996513
$ws.Range("J5").Value = 0.3568846407551645
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.84514
$ws.Range("N5").Value = 3.69028
$ws.Range("O5").Value = 0.0588984363795283
$ws.Range("P5").Value = 0.04005195679609568
$ws.Range("Q5").Value = 7.893096223686666
$ws.Range("R5").Value = 47.35857734212
$ws.Range("S5").Value = 0.02083036165496009
$ws.Range("T5").Value = 0.01429392821271598

$ws.Range("G6").Value = 4.277776333333333
$ws.Range("H6").Value = 12.833329
$ws.Range("I6").Value = 0.3536657835996513
$ws.Range("J6").Value = 0.3568846407551645
$ws.Range("M6").Value = 10.04801166666667
$ws.Range("N6").Value = 30.144035
$ws.Range("O6").Value = 0.3207410689107194
$ws.Range("P6").Value = 0.3271642226280922
$ws.Range("Q6").Value = 42.98314650472388
$ws.Range("R6").Value = 386.848318542515
$ws.Range("S6").Value = 0.1134351414688994
$ws.Range("T6").Value = 0.1167598860605693

$ws.Range("G7").Value = 0.3272805
$ws.Range("H7").Value = 0.6545609999999999
$ws.Range("I7").Value = 0.02705796317293487
$ws.Range("J7").Value = 0.01820281918567982
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.2791963333333333
$ws.Range("N7").Value = 0.837589
$ws.Range("O7").Value = 0.008912184157424861
$ws.Range("P7").Value = 0.009090659364840875
$ws.Range("Q7").Value = 0.09137551557149999
$ws.Range("R7").Value = 0.548253093429
$ws.Range("S7").Value = 0.0002411455507220155
$ws.Range("T7").Value = 0.0001654756286968054

$ws.Range("G8").Value = 0.3272805
$ws.Range("H8").Value = 0.6545609999999999
$ws.Range("I8").Value = 0.02705796317293487
$ws.Range("J8").Value = 0.01820281918567982
$ws.Range("M8").Value = 2.510701
$ws.Range("N8").Value = 7.532103
$ws.Range("O8").Value = 0.08014370894160773
$ws.Range("P8").Value = 0.0817486651256118
$ws.Range("Q8").Value = 0.8217034786305
$ws.Range("R8").Value = 4.930220871783
$ws.Range("S8").Value = 0.002168525525084433
$ws.Range("T8").Value = 0.001488056169952201

$ws.Range("G9").Value = 0.3272805
$ws.Range("H9").Value = 0.6545609999999999
$ws.Range("I9").Value = 0.02705796317293487
$ws.Range("J9").Value = 0.01820281918567982
$ws.Range("M9").Value = 16.644438
$ws.Range("N9").Value = 49.933314
$ws.Range("O9").Value = 0.5313046016107197
$ws.Range("P9").Value = 0.5419444960853593
$ws.Range("Q9").Value = 5.447399990859
$ws.Range("R9").Value = 32.684399945154
$ws.Range("S9").Value = 0.01437602034399369
$ws.Range("T9").Value = 0.00986491767091616

$ws.Range("G10").Value = 0.3272805
$ws.Range("H10").Value = 0.6545609999999999
$ws.Range("I10").Value = 0.02705796317293487
$ws.Range("J10").Value = 0.01820281918567982
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 1.84514
$ws.Range("N10").Value = 3.69028
$ws.Range("O10").Value = 0.0588984363795283
$ws.Range("P10").Value = 0.04005195679609568
$ws.Range("Q10").Value = 0.60387834177
$ws.Range("R10").Value = 2.41551336708
$ws.Range("S10").Value = 0.001593671722500724
$ws.Range("T10").Value = 0.0007290585275919896

$ws.Range("G11").Value = 0.3272805
$ws.Range("H11").Value = 0.6545609999999999
$ws.Range("I11").Value = 0.02705796317293487
$ws.Range("J11").Value = 0.01820281918567982
$ws.Range("M11").Value = 10.04801166666667
$ws.Range("N11").Value = 30.144035
$ws.Range("O11").Value = 0.3207410689107194
$ws.Range("P11").Value = 0.3271642226280922
$ws.Range("Q11").Value = 3.2885182822725
$ws.Range("R11").Value = 19.731109693635
$ws.Range("S11").Value = 0.008678600030634011
$ws.Range("T11").Value = 0.005955311188522661

$ws.Range("G12").Value = 7.490476666666666
$ws.Range("H12").Value = 22.47143
$ws.Range("I12").Value = 0.6192762532274139
$ws.Range("J12").Value = 0.6249125400591558
$ws.Range("K12").Value = 2
$ws.Range("L12").Value = 0.6666666666666666
$ws.Range("M12").Value = 0.2791963333333333
$ws.Range("N12").Value = 0.837589
$ws.Range("O12").Value = 0.008912184157424861
$ws.Range("P12").Value = 0.009090659364840875
$ws.Range("Q12").Value = 2.091313620252222
$ws.Range("R12").Value = 18.82182258227
$ws.Range("S12").Value = 0.005519104013082784
$ws.Range("T12").Value = 0.005680867034495263

$ws.Range("G13").Value = 7.490476666666666
$ws.Range("H13").Value = 22.47143
$ws.Range("I13").Value = 0.6192762532274139
$ws.Range("J13").Value = 0.6249125400591558
$ws.Range("M13").Value = 2.510701
$ws.Range("N13").Value = 7.532103
$ws.Range("O13").Value = 0.08014370894160773
$ws.Range("P13").Value = 0.0817486651256118
$ws.Range("Q13").Value = 18.80634725747667
$ws.Range("R13").Value = 169.25712531729
$ws.Range("S13").Value = 0.04963109579310722
$ws.Range("T13").Value = 0.05108576597009139

$ws.Range("G14").Value = 7.490476666666666
$ws.Range("H14").Value = 22.47143
$ws.Range("I14").Value = 0.6192762532274139
$ws.Range("J14").Value = 0.6249125400591558
$ws.Range("M14").Value = 16.644438
$ws.Range("N14").Value = 49.933314
$ws.Range("O14").Value = 0.5313046016107197
$ws.Range("P14").Value = 0.5419444960853593
$ws.Range("Q14").Value = 124.67477446878
$ws.Range("R14").Value = 1122.07297021902
$ws.Range("S14").Value = 0.3290243230079703
$ws.Range("T14").Value = 0.3386679116197811

$ws.Range("G15").Value = 7.490476666666666
$ws.Range("H15").Value = 22.47143
$ws.Range("I15").Value = 0.6192762532274139
$ws.Range("J15").Value = 0.6249125400591558
$ws.Range("K15").Value = 2
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 1.84514
$ws.Range("N15").Value = 3.69028
$ws.Range("O15").Value = 0.0588984363795283
$ws.Range("P15").Value = 0.04005195679609568
$ws.Range("Q15").Value = 13.82097811673333
$ws.Range("R15").Value = 82.92586870039999
$ws.Range("S15").Value = 0.03647440300206749
$ws.Range("T15").Value = 0.02502897005578772

$ws.Range("G16").Value = 7.490476666666666
$ws.Range("H16").Value = 22.47143
$ws.Range("I16").Value = 0.6192762532274139
$ws.Range("J16").Value = 0.6249125400591558
$ws.Range("M16").Value = 10.04801166666667
$ws.Range("N16").Value = 30.144035
$ws.Range("O16").Value = 0.3207410689107194
$ws.Range("P16").Value = 0.3271642226280922
$ws.Range("Q16").Value = 75.26439693556111
$ws.Range("R16").Value = 677.37957242005
$ws.Range("S16").Value = 0.1986273274111861
$ws.Range("T16").Value = 0.2044490253790003
